$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new data rows (10 and 11) to column A, copying the
# number format/style used by the existing A3:A9 cells.
$ws.Range("A10").Value = 0.11260416666666667
$ws.Range("A11").Value = 0.056365740740740744
$ws.Range("A10:A11").NumberFormat = $ws.Range("A9").NumberFormat

# Extend the SUM formulas so they cover the newly added rows.
$ws.Range("C2").Formula = "=SUM(A2:A11)"
$ws.Range("B3").Formula = "=SUM(A9:A11)"

# Update the active selection to match the author's saved state.
$ws.Range("F3").Select()
